# EI Variable Installments T1 scenarios
# - Expand the "remembered" selection rectangles on the Summary and
#   Repayment schedule sheets.
# - Add a new "O" column of zeros (rows 2-14) to the Repayment schedule
#   sheet, re-using the formatting already applied to the adjacent "N"
#   column, and shift the previously-empty formatting-only cell that sat
#   at the end of the header row ("O2") one column to the right ("P2").

$wb = $excel.ActiveWorkbook

# Remember which sheet is active/selected so we can restore it at the
# end - selecting ranges on other sheets below would otherwise change
# the workbook's active tab.
$originalActiveSheet = $wb.ActiveSheet

$wsSummary = $wb.Worksheets.Item("Summary")
$wsRepay   = $wb.Worksheets.Item("Repayment schedule")

# --- Repayment schedule: insert the new "O" column values ------------
# Row 2 only carries formatting (no value) on O2; the new P2 cell is the
# same - formatting only, copied from O2.
$wsRepay.Range("O2").Copy()
$wsRepay.Range("P2").PasteSpecial(-4122)  # xlPasteFormats
$wsRepay.Range("P2").ClearContents()

# Rows 3-14: new "O" cell, value 0, formatted like the existing "N"
# cell on the same row.
for ($r = 3; $r -le 14; $r++) {
    $wsRepay.Range("N$r").Copy()
    $wsRepay.Range("O$r").PasteSpecial(-4122)  # xlPasteFormats
    $wsRepay.Range("O$r").Value = 0
}

$excel.CutCopyMode = 0

# --- Selections (sheetView "remembered" selection) --------------------
$wsSummary.Range("A7:XFD16").Select()
$wsRepay.Range("A15:XFD15").Select()

# Restore the sheet that was active before we touched other sheets.
$originalActiveSheet.Activate()
